# Commit: "Analytical solution succeed, but volumetric strain prediction is not so good."
#
# The author reduced the AP column's constant coefficient (used across rows 3-18 of
# Sheet1) from ~28139.5986 to 22197.0184441115, and switched that column's number
# format to scientific notation (0.00E+00) to make the smaller coefficient legible.
# That ripples through the AU/AV/AW/AX (and dependent AO/AT/AQ/AN/AJ/...) formula
# columns via recalculation. The author also added a new summary row (row 23) with
# a MIN() of the W column, and left the selection on cell I16 after scrolling back
# to the left edge of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the AP coefficient for rows 3 through 18 and apply scientific notation
#    number formatting (this creates/uses the new cellXfs style entry: numFmtId=11,
#    fillId=2, applyNumberFormat=1, applyFill=1).
$apRange = $ws.Range("AP3:AP18")
$apRange.Value = 22197.018444111502
$apRange.NumberFormat = "0.00E+00"

# 2) Add the new MIN() summary row right below the data (row 23).
$ws.Range("W23").Formula = "=MIN(W2:W22)"

# 3) Leave the sheet scrolled back to the left and the selection on I16 (matches
#    the author's final view state after finishing the edits).
$ws.Range("I16").Select()

$excel.CalculateFull()
